$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.116748213768005
$ws.Range("B1").Value = 1.735088706016541
$ws.Range("C1").Value = 4.522050380706787
$ws.Range("D1").Value = 0.3457743525505066
$ws.Range("E1").Value = 0.3952162861824036
